# Fruta / hortaliza, semanal
# Insert a new weekly record at row 170, pushing the existing rows
# (170-247) down by one (170->171, ..., 247->248).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("170:170").Insert()

$ws.Range("A170").Value = 5
$ws.Range("B170").Value = "Macroferia Regional de Talca"
$ws.Range("C170").Value = "Maule"
$ws.Range("D170").Value = 44726
$ws.Range("E170").Value = 7
$ws.Range("F170").Value = 100112024
$ws.Range("G170").Value = "Choclo"
$ws.Range("H170").Value = "Dulce o Americano"
$ws.Range("I170").Value = "Primera"
$ws.Range("J170").Value = 150
$ws.Range("K170").Value = 32000
$ws.Range("L170").Value = 32000
$ws.Range("M170").Value = 32000
$ws.Range("N170").Value = "$/malla 60 unidades"
$ws.Range("O170").Value = "Provincia del Elquí"
$ws.Range("P170").Value = 533
$ws.Range("Q170").Value = 60
$ws.Range("R170").Value = "Hortaliza"
